$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: fill in G21/H21 (existing style retained) and add J21 (style copied from J7, which
#     already carries the "thick border" cell style used for this new cell) ---
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 5

$ws.Range("J7").Copy()
$ws.Range("J21").PasteSpecial(-4122)
$ws.Range("J21").Value = 5

# --- Row 22: add I22 (style copied from I5, the green-fill/thick-border style) ---
$ws.Range("I5").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$ws.Range("I22").Value = 5

# --- Row 23: remove C23:E23 entirely, bump F23 to the thick-border style (copied from I4),
#     fill G23/H23, add J23 (style copied from J9) and the N23 remark (shared string reuse) ---
$ws.Range("C23:E23").Clear()

$ws.Range("I4").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("F23").Value = 5

$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 5

$ws.Range("J9").Copy()
$ws.Range("J23").PasteSpecial(-4122)
$ws.Range("J23").Value = 5

$ws.Range("N23").Value = "четверку очень нада"

# --- Conditional formatting: new 3-colour scale on E23, inserted ahead of the existing rules
#     (K4:K30 and E30 shift down in priority, matching Excel's "new rule goes first" behaviour) ---
$cs = $ws.Range("E23").FormatConditions.AddColorScale(3)
$cs.SetFirstPriority()

# --- Selection / scroll state: land on I22 with the frozen pane scrolled back to the top ---
$ws.Range("I22").Select()
